$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.750.63"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "3.500.95"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'599.45"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").Value = "'182.42"
$ws.Range("E6").Value = "  +5.41%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.596"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("D9").Value = "'0.140"
$ws.Range("E9").Value = "  +5.03%  "
$ws.Range("D10").Value = "'7.12"
$ws.Range("E10").Value = "  -2.21%  "
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "4.110.85"
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("D13").Value = "'32.44"
$ws.Range("E13").Value = "  +12.75%  "
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "67.748.77"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000182"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").Value = "3.502.08"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("D19").Value = "'14.72"
$ws.Range("E19").Value = "  +3.10%  "
$ws.Range("D20").Value = "'396.43"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("E21").Value = "  +1.48%  "
$ws.Range("D22").Value = "'73.51"
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").Value = "'0.547"
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "'0.0000126"
$ws.Range("E25").Value = "  +3.01%  "
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "'5.70"
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").Value = "'10.49"
$ws.Range("E27").Value = "  +2.55%  "
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").Value = "'1.47"
$ws.Range("E31").Value = "  +0.96%  "
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("D33").Value = "'23.98"
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").Value = "'7.44"
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").Value = "'1.67"
$ws.Range("E36").Value = "  +1.58%  "
$ws.Range("D37").Value = "'164.55"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("E38").Value = "  +2.38%  "
$ws.Range("D39").Value = "'0.875"
$ws.Range("E39").Value = "  -1.93%  "
$ws.Range("E40").Value = "  +3.18%  "
$ws.Range("D41").Value = "'4.74"
$ws.Range("E41").Value = "  +1.04%  "
$ws.Range("D42").Value = "'27.86"
$ws.Range("E42").Value = "  +1.95%  "
$ws.Range("D43").Value = "'2.69"
$ws.Range("E43").Value = "  +3.22%  "
$ws.Range("D44").Value = "'26.67"
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("D45").Value = "'0.0738"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("D46").Value = "2.814.16"
$ws.Range("E46").Value = "  +0.67%  "
$ws.Range("D47").Value = "'42.40"
$ws.Range("E47").Value = "  -1.11%  "
$ws.Range("D48").Value = "'0.0306"
$ws.Range("E48").Value = "  -0.48%  "
$ws.Range("D49").Value = "'345.08"
$ws.Range("E49").Value = "  +1.18%  "
$ws.Range("E50").Value = "  -0.91%  "
$ws.Range("D51").Value = "'33.71"
$ws.Range("E51").Value = "  +0.42%  "
